# db_config.xlsx edit
# - F2 ("report_type varchar (225)") becomes E2 with updated text "report_testing int"
# - F2 is left as a blank (but touched/formatted) cell
# - Two new blank rows are added below (A3, A4), mirroring A2's formatting
# - Selection moves to E2, view scrolled back to the top-left (A1)

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Move / retype the value that used to live in F2 into E2.
$ws.Range("E2").Value = "report_testing int"

# F2 no longer holds the text, but remains a live (formatted) cell in the sheet.
$ws.Range("F2").ClearContents()
$ws.Range("F2").NumberFormat = "General"

# New rows 3 and 4 appear with a formatted-but-empty cell in column A,
# matching the look of A2.
$ws.Range("A3").NumberFormat = "General"
$ws.Range("A4").NumberFormat = "General"

# Update the active selection / view to match (active cell E2, view scrolled to A1).
[void]$ws.Range("E2").Select()
